$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.571.45'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '3.193.63'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''596.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.53%  '
$ws.Range('D6').Value = '''154.65'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.06%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.563'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +7.32%  '
$ws.Range('D9').Value = '3.179.97'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('D11').Value = '''5.93'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.88%  '
$ws.Range('D12').Value = '''0.520'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.04%  '
$ws.Range('D13').Value = '''0.0000269'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('D14').Value = '''39.29'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +5.99%  '
$ws.Range('D15').Value = '3.716.46'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.559.18'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''7.50'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +5.41%  '
$ws.Range('D18').Value = '3.192.21'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = '''517.80'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.95%  '
$ws.Range('E21').Value = '  +3.87%  '
$ws.Range('D22').Value = '''0.740'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.08%  '
$ws.Range('D23').Value = '''8.09'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.32%  '
$ws.Range('D24').Value = '''14.97'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = '''86.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '''9.29'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +4.92%  '
$ws.Range('E28').Value = '  +3.92%  '
$ws.Range('E29').Value = '  +8.87%  '
$ws.Range('D30').Value = '''7.11'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +15.54%  '
$ws.Range('E31').Value = '  +5.39%  '
$ws.Range('D32').Value = '''28.39'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').Value = '''1.23'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.32%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = '''6.54'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.60%  '
$ws.Range('D36').Value = '''517.01'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +8.24%  '
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('E39').Value = '  +2.86%  '
$ws.Range('D40').Value = '''0.129'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +12.40%  '
$ws.Range('D41').Value = '''8.92'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('D42').Value = '''2.91'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').Value = '''0.303'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +8.01%  '
$ws.Range('D44').Value = '0.0₃0672'
$ws.Range('E44').Value = '  +16.23%  '
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('D46').Value = '2.908.15'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('D47').Value = '''28.88'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.36%  '
$ws.Range('E48').Value = '  +4.17%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '''0.999'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '''2.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.29%  '
$ws.Range('E51').Value = '  +9.89%  '
